# Fruta / hortaliza, semanal
# Update weekly market data: dates and associated volume/price figures were
# reshuffled across the existing rows for the "Espárragos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry lists the new values for D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) per row.
$rowsData = @{
    2  = @{ D = 44179; J = 200; K = 1600; L = 1600; M = 1600; P = 1600 }
    3  = @{ D = 44159; J = 600; K = 1600; L = 1700; M = 1650; P = 1650 }
    4  = @{ D = 44169; J = 600; K = 1600; L = 1600; M = 1600; P = 1600 }
    5  = @{ D = 44161; J = 300; K = 1700; L = 1700; M = 1700; P = 1700 }
    6  = @{ D = 44165; J = 300; K = 1600; L = 1600; M = 1600; P = 1600 }
    7  = @{ D = 44166; J = 500; K = 1600; L = 1600; M = 1600; P = 1600 }
    8  = @{ D = 44168; J = 200; K = 1600; L = 1600; M = 1600; P = 1600 }
    9  = @{ D = 44162; J = 700; K = 1600; L = 1600; M = 1600; P = 1600 }
    10 = @{ D = 44176; J = 700; K = 1600; L = 1600; M = 1600; P = 1600 }
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
